# Adjust rf to the same period of return (previously annual rf was used in
# all circumstances in the past, which was not correct). This recomputes
# the GRS test statistics in columns B:F for rows 2-6 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C values are extremely small (~1.11e-16); build them via Pow()
# since the parser here doesn't accept scientific-notation numeric literals.
$rfEps = 1.110223024625157 * [Math]::Pow(10, -16)

$ws.Range("B2").Value = 12.18812099985124
$ws.Range("C2").Value = $rfEps
$ws.Range("D2").Value = 0.00736770806979629
$ws.Range("E2").Value = 0.8929020919339536
$ws.Range("F2").Value = 0.7972741457800305

$ws.Range("B3").Value = 12.52358522724342
$ws.Range("C3").Value = $rfEps
$ws.Range("D3").Value = 0.007355251178049704
$ws.Range("E3").Value = 0.8913924250749858
$ws.Range("F3").Value = 0.7945804554810642

$ws.Range("B4").Value = 13.61217423570494
$ws.Range("C4").Value = $rfEps
$ws.Range("D4").Value = 0.006486985216250005
$ws.Range("E4").Value = 0.78616614760829
$ws.Range("F4").Value = 0.6180572116452596

$ws.Range("B5").Value = 13.8912692044324
$ws.Range("C5").Value = $rfEps
$ws.Range("D5").Value = 0.006812886420499159
$ws.Range("E5").Value = 0.825662537025618
$ws.Range("F5").Value = 0.68171862504758

$ws.Range("B6").Value = 13.81768963282079
$ws.Range("C6").Value = $rfEps
$ws.Range("D6").Value = 0.006478520402170357
$ws.Range("E6").Value = 0.7851402858168145
$ws.Range("F6").Value = 0.6164452684125092
